$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "old"/"new" AHB version headers to the actual format versions
# (FV2410 = old/reference version, FV2504 = new version being compared).
$headerRange = $ws.Range("A1:U1")
[void]$headerRange.Replace("_old", "_FV2410")
[void]$headerRange.Replace("_new", "_FV2504")

# Freeze the header row.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a proper Excel Table so it can be filtered/sorted.
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
